# Daily update at 8 AM UTC
# Appends the new day's win counts as the next row of the log, then fixes
# up the date number-formatting so that only the newest (last) row keeps
# the "date-only" display format while the row that used to be last goes
# back to the normal "date + time" format used by every other data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row/column of the existing data table.
$usedRange = $ws.UsedRange
$lastRow = $ws.Cells.Item($usedRange.Rows.Count, 1).Row
$newRow = $lastRow + 1

# Today's day serial is one more than the previous day's. Value2() returns
# the raw underlying serial number rather than a formatted display string.
$prevDay = $ws.Cells.Item($lastRow, 1).Value2()
$newDay = $prevDay + 1

# New win totals for the day.
$chaseWins = 380
$bryceWins = 384
$zachWins = 386

$ws.Cells.Item($newRow, 1).Value = $newDay
$ws.Cells.Item($newRow, 2).Value = $chaseWins
$ws.Cells.Item($newRow, 3).Value = $bryceWins
$ws.Cells.Item($newRow, 4).Value = $zachWins

# The previously-last row reverts to the standard date+time format, and the
# newly added last row takes on the "latest row" date-only format.
$ws.Cells.Item($lastRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD"
